$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2000-2009 (original rows 2 through 11).
# This shifts the 2010/2011/2012 rows (originally 12,13,14) up to rows 2,3,4.
$ws.Range("A2:F11").EntireRow.Delete()
